$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so values such as
# "29.122.16" or "1.000" are not reinterpreted as numbers/dates.
$dCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.122.16"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.841.60"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "241.28"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "0.6861"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.3022"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").Value = "23.15"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").Value = "1.834.11"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "5.060"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "0.6834"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "87.57"
$ws.Range("E15").Value = "  -6.20%  "
$ws.Range("E16").Value = "  -7.09%  "
$ws.Range("D17").Value = "29.116.38"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "0.000008163"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").Value = "2.080.06"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "227.80"
$ws.Range("E20").Value = "  -5.76%  "
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "7.391"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "0.1457"
$ws.Range("E25").Value = "  -3.58%  "
$ws.Range("D26").Value = "160.30"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "8.756"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").Value = "1.513"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "4.265"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "4.148"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "1.194"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "0.05190"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").Value = "0.7648"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").Value = "1.850"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "1.317.17"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").Value = "2.727"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "0.9357"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").Value = "105.10"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "5.769"
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.980.56"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.5193"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "64.93"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "9.502"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "0.05942"
$ws.Range("E51").Value = "  +0.95%  "
